$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts existing A:D to B:E, carrying
# their column widths and formatting along)
$ws.Columns.Item(1).Insert()

# Insert a new row before row 1 (shifts existing rows 1:14 down to 2:15,
# carrying their formatting along)
$ws.Rows.Item(1).Insert()

# New header row (row 1) across the data columns B:E
$ws.Cells.Item(1, 2).Value = "Valid"
$ws.Cells.Item(1, 3).Value = "T"
$ws.Cells.Item(1, 4).Value = "Z"
$ws.Cells.Item(1, 5).Value = "p-value"

# New row labels in column A for the 14 data rows (now rows 2-15)
$labels = @(
    "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)",
    "MaintainabilityIndex & MaintainabilityIndex",
    "NbOperands & NbOperands",
    "NbOperands & EffortToImplement",
    "NbUniqueOperators & NbUniqueOperators",
    "ProgramLength & EffortToImplement",
    "VocabularySize & VocabularySize",
    "ProgramVolume & ProgramVolume",
    "DifficultyLevel & DifficultyLevel",
    "ProgramLevel & ProgramLevel",
    "EffortToImplement & NbOperands",
    "EffortToImplement & ProgramLength",
    "EffortToImplement & EffortToImplement",
    "TimeToImplement & TimeToImplement"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# New column A is wide (holds the long dataset-pair labels); the other
# columns (B:E) already kept their original widths via the column insert.
$ws.Columns.Item(1).ColumnWidth = 53.67
